$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-21 Friday" "2025-02-22 Saturday"

Replace-Text "13×57=" "66×73="
Replace-Text "25×80=" "65×28="
Replace-Text "22×20=" "66×29="
Replace-Text "65×86=" "97×76="
Replace-Text "97×85=" "89×40="

Replace-Text "77×92=" "52×80="
Replace-Text "32×13=" "64×68="
Replace-Text "65×21=" "54×25="
Replace-Text "63×22=" "53×46="
Replace-Text "74×94=" "35×11="

Replace-Text "28×23=" "52×62="
Replace-Text "24×14=" "39×60="
Replace-Text "29×50=" "37×37="
Replace-Text "67×63=" "78×91="
Replace-Text "79×62=" "76×36="

Replace-Text "38×17=" "89×18="
Replace-Text "62×78=" "47×97="
Replace-Text "74×22=" "78×75="
Replace-Text "24×62=" "46×83="
Replace-Text "38×80=" "41×17="

Replace-Text "13×44=" "82×30="
Replace-Text "85×77=" "33×93="
Replace-Text "70×80=" "65×32="
Replace-Text "72×78=" "36×95="
Replace-Text "74×65=" "89×73="
